# Refresh the cryptocurrency price / 1h-volume figures to match the
# latest scrape (GitHub Actions run). A handful of rows also swapped
# rank position (PaxDollar <-> Quant).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.460.22"
$ws.Range("E2").Value = "'  +2.06%  "
$ws.Range("D3").Value = "'1.838.42"
$ws.Range("D5").Value = "'314.33"
$ws.Range("E5").Value = "'  +1.63%  "
$ws.Range("E6").Value = "'  +0.93%  "
$ws.Range("D7").Value = "'0.4737"
$ws.Range("E7").Value = "'  +1.66%  "
$ws.Range("D8").Value = "'0.3692"
$ws.Range("E8").Value = "'  +0.73%  "
$ws.Range("D9").Value = "'0.07459"
$ws.Range("E9").Value = "'  +1.51%  "
$ws.Range("D10").Value = "'0.8855"
$ws.Range("E10").Value = "'  +2.02%  "
$ws.Range("E11").Value = "'  +0.63%  "
$ws.Range("D12").Value = "'1.912.72"
$ws.Range("E12").Value = "'  +3.40%  "
$ws.Range("D13").Value = "'0.07326"
$ws.Range("E13").Value = "'  +3.08%  "
$ws.Range("D14").Value = "'5.451"
$ws.Range("E14").Value = "'  +1.29%  "
$ws.Range("D15").Value = "'93.29"
$ws.Range("E15").Value = "'  +1.78%  "
$ws.Range("D16").Value = "'6.584"
$ws.Range("E16").Value = "'  +1.05%  "
$ws.Range("D17").Value = "'1.010"
$ws.Range("E17").Value = "'  +0.70%  "
$ws.Range("D18").Value = "'0.000008816"
$ws.Range("E18").Value = "'  +1.31%  "
$ws.Range("E19").Value = "'  +0.93%  "
$ws.Range("E20").Value = "'  +1.11%  "
$ws.Range("D21").Value = "'27.487.66"
$ws.Range("E21").Value = "'  +2.03%  "
$ws.Range("D22").Value = "'5.325"
$ws.Range("E22").Value = "'  +0.47%  "
$ws.Range("E23").Value = "'  +0.48%  "
$ws.Range("D24").Value = "'2.140.99"
$ws.Range("E24").Value = "'  +2.80%  "
$ws.Range("D25").Value = "'1.908"
$ws.Range("E25").Value = "'  +0.79%  "
$ws.Range("D26").Value = "'152.16"
$ws.Range("E26").Value = "'  +0.76%  "
$ws.Range("D27").Value = "'18.64"
$ws.Range("E27").Value = "'  +1.86%  "
$ws.Range("D28").Value = "'2.148"
$ws.Range("E28").Value = "'  +0.40%  "
$ws.Range("D29").Value = "'5.253"
$ws.Range("E29").Value = "'  -0.13%  "
$ws.Range("D31").Value = "'0.08998"
$ws.Range("E31").Value = "'  +1.12%  "
$ws.Range("D32").Value = "'0.7569"
$ws.Range("E32").Value = "'  +0.22%  "
$ws.Range("D33").Value = "'1.182"
$ws.Range("E33").Value = "'  +2.27%  "
$ws.Range("D34").Value = "'4.561"
$ws.Range("E34").Value = "'  +1.62%  "
$ws.Range("D35").Value = "'2.950"
$ws.Range("E35").Value = "'  +1.29%  "
$ws.Range("D36").Value = "'1.012"
$ws.Range("E36").Value = "'  +1.10%  "
$ws.Range("D37").Value = "'1.105"
$ws.Range("E37").Value = "'  +1.90%  "
$ws.Range("D38").Value = "'0.05330"
$ws.Range("E38").Value = "'  +1.04%  "
$ws.Range("E39").Value = "'  +0.42%  "
$ws.Range("D40").Value = "'2.996"
$ws.Range("E40").Value = "'  +0.49%  "
$ws.Range("D41").Value = "'7.327"
$ws.Range("E41").Value = "'  +1.05%  "
$ws.Range("D42").Value = "'2.404"
$ws.Range("D43").Value = "'0.5333"
$ws.Range("E43").Value = "'  +0.60%  "
$ws.Range("D44").Value = "'0.1660"
$ws.Range("E44").Value = "'  +0.45%  "
$ws.Range("D45").Value = "'8.515"
$ws.Range("E45").Value = "'  +1.04%  "
$ws.Range("D46").Value = "'0.4911"
$ws.Range("E46").Value = "'  +0.82%  "
$ws.Range("E47").Value = "'  +1.28%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "'1.012"
$ws.Range("E48").Value = "'  +1.08%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'104.94"
$ws.Range("E49").Value = "'  +1.63%  "
$ws.Range("E50").Value = "'  +1.15%  "
$ws.Range("D51").Value = "'0.06316"
$ws.Range("E51").Value = "'  +0.40%  "
